$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = '[Ian%Leonard-Lorant%NULL%0, Xavier%Delabranche%NULL%3, Xavier%Delabranche%NULL%0, François%Severac%NULL%3, François%Severac%NULL%0, Julie%Helms%NULL%0, Julie%Helms%NULL%0, Coralie%Pauzet%NULL%1, Olivier%Collange%NULL%1, Françis%Schneider%NULL%1, Aissam%Labani%NULL%2, Aissam%Labani%NULL%0, Pascal%Bilbault%NULL%1, Sébastien%Moliere%NULL%2, Sébastien%Moliere%NULL%0, Pierre%Leyendecker%NULL%1, Catherine%Roy%NULL%2, Catherine%Roy%NULL%0, Mickaël%Ohana%mickael.ohana@gmail.com%3, Mickaël%Ohana%mickael.ohana@gmail.com%0]'
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = 'Radiological Society of North America'

# Row 3
$ws.Range("E3").Value = '[Franck%Grillet%f1grillet@chu-besancon.fr%0, Julien%Behr%NULL%2, Julien%Behr%NULL%0, Paul%Calame%NULL%1, Sébastien%Aubry%NULL%1, Eric%Delabrousse%NULL%2, Eric%Delabrousse%NULL%0]'
$ws.Range("I3").Value = ''
$ws.Range("J3").Value = 'Radiological Society of North America'

# Row 4
$ws.Range("D4").Value = 'Objectives
id="Par1">To evaluate the prevalence of acute pulmonary embolism (APE) in non-hospitalized COVID-19 patients referred to CT pulmonary angiography (CTPA) by the emergency department.


Methods
id="Par2">From March 14 to April 6, 2020, 72 non-hospitalized patients referred by the emergency department to CTPA for COVID-19 pneumonia were retrospectively identified.

 Relevant clinical and laboratory data and CT scan findings were collected for each patient.

 CTPA scans were reviewed by two radiologists to determinate the presence or absence of APE.

 Clinical classification, lung involvement of COVID-19 pneumonia, and CT total severity score were compared between APE group and non-APE group.


Results
id="Par3">APE was identified in 13 (18%) CTPA scans.

 The mean age and D-dimer of patients from the APE group were higher in comparison with those from the non-APE group (74.4 vs.

 59.6 years, p = 0.008, and 7.29 vs.

 3.29 μg/ml, p = 0.011).

 There was no significant difference between APE and non-APE groups concerning clinical type, COVID-19 pneumonia lung lesions (ground-glass opacity: 85% vs.

 97%; consolidation: 69% vs.

 68%; crazy paving: 38% vs.

 37%; linear reticulation: 69% vs.

 78%), CT severity score (6.3 vs.

 7.1, p = 0.365), quality of CTPA (1.8 vs.

 2.0, p = 0.518), and pleural effusion (38% vs.

 19%, p = 0.146).


Conclusions
id="Par4">Non-hospitalized patients with COVID-19 pneumonia referred to CT scan by the emergency departments are at risk of APE.

 The presence of APE was not limited to severe or critical clinical type of COVID-19 pneumonia.


Key Points
id="Par5">
• Acute pulmonary embolism was found in 18% of non-hospitalized COVID-19 patients referred by the emergency department to CTPA.

 Two (15%) patients had main, four (30%) lobar, and seven (55%) segmental acute pulmonary embolism.



'
$ws.Range("E4").Value = '[Alban%Gervaise%alban.gervaise@hotmail.fr%0, Caroline%Bouzad%NULL%2, Caroline%Bouzad%NULL%0, Evelyne%Peroux%NULL%1, Carole%Helissey%NULL%1]'
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = 'Springer Berlin Heidelberg'

# Row 5
$ws.Range("D5").Value = 'Patients with coronavirus disease 2019 (COVID-19) have elevated D-dimer levels.
 Early reports describe high venous thromboembolism (VTE) and disseminated intravascular coagulation (DIC) rates, but data are limited.
 This multicenter retrospective study describes the rate and severity of hemostatic and thrombotic complications of 400 hospital-admitted COVID-19 patients (144 critically ill) primarily receiving standard-dose prophylactic anticoagulation.
 Coagulation and inflammatory parameters were compared between patients with and without coagulation-associated complications.
 Multivariable logistic models examined the utility of these markers in predicting coagulation-associated complications, critical illness, and death.
 The radiographically confirmed VTE rate was 4.8% (95% confidence interval [CI], 2.9-7.3), and the overall thrombotic complication rate was 9.5% (95% CI, 6.8-12.8).
 The overall and major bleeding rates were 4.8% (95% CI, 2.9-7.3) and 2.3% (95% CI, 1.0-4.2), respectively.
 In the critically ill, radiographically confirmed VTE and major bleeding rates were 7.6% (95% CI, 3.9-13.3) and 5.6% (95% CI, 2.4-10.7), respectively.
 Elevated D-dimer at initial presentation was predictive of coagulation-associated complications during hospitalization (D-dimer &gt;2500 ng/mL, adjusted odds ratio [OR] for thrombosis, 6.79 [95% CI, 2.39-19.30]; adjusted OR for bleeding, 3.56 [95% CI, 1.01-12.66]), critical illness, and death.
 Additional markers at initial presentation predictive of thrombosis during hospitalization included platelet count &gt;450 × 109/L (adjusted OR, 3.56 [95% CI, 1.27-9.97]), C-reactive protein (CRP) &gt;100 mg/L (adjusted OR, 2.71 [95% CI, 1.26-5.86]), and erythrocyte sedimentation rate (ESR) &gt;40 mm/h (adjusted OR, 2.64 [95% CI, 1.07-6.51]).
 ESR, CRP, fibrinogen, ferritin, and procalcitonin were higher in patients with thrombotic complications than in those without.
 DIC, clinically relevant thrombocytopenia, and reduced fibrinogen were rare and were associated with significant bleeding manifestations.
 Given the observed bleeding rates, randomized trials are needed to determine any potential benefit of intensified anticoagulant prophylaxis in COVID-19 patients.
'
$ws.Range("E5").Value = '[Hanny%Al-Samkari%NULL%0, Rebecca S.%Karp Leaf%NULL%1, Walter H.%Dzik%NULL%1, Jonathan C.T.%Carlson%NULL%1, Annemarie E.%Fogerty%NULL%1, Anem%Waheed%NULL%1, Katayoon%Goodarzi%NULL%1, Pavan K.%Bendapudi%NULL%1, Larissa%Bornikova%NULL%1, Shruti%Gupta%NULL%1, David E.%Leaf%NULL%1, David J.%Kuter%NULL%1, Rachel P.%Rosovsky%NULL%1]'
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = 'American Society of Hematology'

# Row 6
$ws.Range("C6").Value = 'Unknown Title'
$ws.Range("E6").Value = '[]'
$ws.Range("F6").Value = 'not found'
$ws.Range("G6").Value = 'N/A'
$ws.Range("H6").Value = '1970-01-01'
$ws.Range("J6").Value = ''

# Row 7
$ws.Range("C7").Value = 'Unknown Title'
$ws.Range("D7").Value = 'Unknown Abstract'
$ws.Range("E7").Value = '[]'
$ws.Range("F7").Value = 'not found'
$ws.Range("G7").Value = 'N/A'
$ws.Range("H7").Value = '1970-01-01'
$ws.Range("J7").Value = ''

# Row 8
$ws.Range("E8").Value = '[Jean‐François%Llitjos%NULL%0, Maxime%Leclerc%NULL%0, Camille%Chochois%NULL%1, Jean‐Michel%Monsallier%NULL%1, Michel%Ramakers%NULL%0, Malika%Auvray%NULL%1, Karim%Merouani%NULL%1]'
$ws.Range("I8").Value = ''
$ws.Range("J8").Value = 'International Society on Thrombosis and Haemostasis. Published by Elsevier Inc.'

# Row 9
$ws.Range("D9").Value = 'Purpose
id="Par1">Little evidence of increased thrombotic risk is available in COVID-19 patients.

 Our purpose was to assess thrombotic risk in severe forms of SARS-CoV-2 infection.


Methods
id="Par2">All patients referred to 4 intensive care units (ICUs) from two centers of a French tertiary hospital for acute respiratory distress syndrome (ARDS) due to COVID-19 between March 3rd and 31st 2020 were included.

 Medical history, symptoms, biological data and imaging were prospectively collected.

 Propensity score matching was performed to analyze the occurrence of thromboembolic events between non-COVID-19 ARDS and COVID-19 ARDS patients.


Results
id="Par3">150 COVID-19 patients were included (122 men, median age 63 [53; 71] years, SAPSII 49 [37; 64] points).

 Sixty-four clinically relevant thrombotic complications were diagnosed in 150 patients, mainly pulmonary embolisms (16.7%).

 28/29 patients (96.6%) receiving continuous renal replacement therapy experienced circuit clotting.

 Three thrombotic occlusions (in 2 patients) of centrifugal pump occurred in 12 patients (8%) supported by ECMO.

 Most patients (&amp;gt; 95%) had elevated D-dimer and fibrinogen.

 No patient developed disseminated intravascular coagulation.

 Von Willebrand (vWF) activity, vWF antigen and FVIII were considerably increased, and 50/57 tested patients (87.7%) had positive lupus anticoagulant.

 Comparison with non-COVID-19 ARDS patients (n = 145) confirmed that COVID-19 ARDS patients (n = 77) developed significantly more thrombotic complications, mainly pulmonary embolisms (11.7 vs.

 2.1%, p &amp;lt; 0.008).

 Coagulation parameters significantly differed between the two groups.


Conclusion
id="Par4">Despite anticoagulation, a high number of patients with ARDS secondary to COVID-19 developed life-threatening thrombotic complications.

 Higher anticoagulation targets than in usual critically ill patients should therefore probably be suggested.


Electronic supplementary material
The online version of this article (10.1007/s00134-020-06062-x) contains supplementary material, which is available to authorized users.


'
$ws.Range("E9").Value = '[Julie%Helms%NULL%0, Charles%Tacquard%NULL%1, François%Severac%NULL%0, Ian%Leonard-Lorant%NULL%0, Mickaël%Ohana%NULL%0, Xavier%Delabranche%NULL%0, Hamid%Merdji%NULL%1, Raphaël%Clere-Jehl%NULL%1, Malika%Schenck%NULL%1, Florence%Fagot Gandet%NULL%1, Samira%Fafi-Kremer%NULL%0, Vincent%Castelain%NULL%0, Francis%Schneider%NULL%0, Lélia%Grunebaum%NULL%1, Eduardo%Anglés-Cano%NULL%1, Laurent%Sattler%NULL%1, Paul-Michel%Mertes%NULL%0, Ferhat%Meziani%Ferhat.Meziani@chru-strasbourg.fr%1, NULL%NULL%NULL%0]'
$ws.Range("I9").Value = ''
$ws.Range("J9").Value = 'Springer Berlin Heidelberg'

# Row 10
$ws.Range("E10").Value = '[P.%Demelo-Rodríguez%NULL%0, E.%Cervilla-Muñoz%NULL%1, L.%Ordieres-Ortega%NULL%1, A.%Parra-Virto%NULL%1, M.%Toledano-Macías%NULL%1, N.%Toledo-Samaniego%NULL%1, A.%García-García%NULL%1, I.%García-Fernández-Bravo%NULL%1, Z.%Ji%NULL%1, J.%de-Miguel-Diez%NULL%1, L.A.%Álvarez-Sala-Walther%NULL%1, J.%del-Toro-Cervera%NULL%1, F.%Galeano-Valle%paco.galeano.valle@gmail.com%1]'
$ws.Range("I10").Value = ''
$ws.Range("J10").Value = 'Elsevier Ltd.'

# Row 11
$ws.Range("E11").Value = '[Songping%Cui%NULL%0, Shuo%Chen%NULL%1, Xiunan%Li%NULL%1, Shi%Liu%NULL%1, Feng%Wang%NULL%0]'
$ws.Range("I11").Value = ''
$ws.Range("J11").Value = 'International Society on Thrombosis and Haemostasis. Published by Elsevier Inc.'

# Row 12
$ws.Range("E12").Value = '[Saskia%Middeldorp%NULL%0, Michiel%Coppens%NULL%1, Thijs F.%van Haaps%NULL%1, Merijn%Foppen%NULL%1, Alexander P.%Vlaar%NULL%1, Marcella C.A.%Müller%NULL%1, Catherine C.S.%Bouman%NULL%1, Ludo F.M.%Beenen%NULL%1, Ruud S.%Kootte%NULL%1, Jarom%Heijmans%NULL%1, Loek P.%Smits%NULL%1, Peter I.%Bonta%NULL%1, Nick%van Es%NULL%1]'
$ws.Range("I12").Value = ''
$ws.Range("J12").Value = 'The Authors. Journal of Thrombosis and Haemostasis published by ELSEVIER INC. on behalf of International Society on Thrombosis and Haemostasis'

# Row 13
$ws.Range("C13").Value = 'Unknown Title'
$ws.Range("E13").Value = '[]'
$ws.Range("F13").Value = 'not found'
$ws.Range("G13").Value = 'N/A'
$ws.Range("H13").Value = '1970-01-01'
$ws.Range("J13").Value = ''

# Row 14
$ws.Range("D14").Value = 'Acute respiratory distress syndrome development in patients with coronavirus disease 2019 (COVID-19) pneumonia is associated with a high mortality rate and is the main cause of death in patients with severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) infection [1].
 Myocardial injury has also been reported to be significantly associated with fatal outcome, with a 37% mortality rate in patients without prior cardiovascular disease but elevated troponin levels [2].
 A D-dimer level of &gt;1 μg·mL−1 has been clearly identified as a risk factor for poor outcome in SARS-Cov-2 infection [3], with recent reports highlighting a high incidence of thrombotic events in intensive care unit (ICU) patients [4].
 A normal D-dimer level allows the safe exclusion of pulmonary embolism (PE) in outpatients with a low or intermediate clinical probability of PE, but there is no recommendation to use D-dimer as a positive marker of thrombosis because of lack of specificity.
'
$ws.Range("E14").Value = '[Florian%Bompard%NULL%0, Hippolyte%Monnier%NULL%1, Ines%Saab%NULL%1, Mickael%Tordjman%NULL%1, Hendy%Abdoul%NULL%1, Laure%Fournier%NULL%1, Olivier%Sanchez%NULL%1, Christine%Lorut%NULL%1, Guillaume%Chassagnon%NULL%1, Marie-Pierre%Revel%NULL%1]'
$ws.Range("I14").Value = ''
$ws.Range("J14").Value = 'European Respiratory Society'

# Row 15
$ws.Range("E15").Value = '[Marco%Cattaneo%NULL%0, Elena M.%Bertinato%NULL%1, Simone%Birocchi%NULL%1, Carolina%Brizio%NULL%1, Daniele%Malavolta%NULL%1, Marco%Manzoni%NULL%1, Gesualdo%Muscarella%NULL%1, Michela%Orlandi%NULL%1]'
$ws.Range("I15").Value = ''
$ws.Range("J15").Value = 'Georg Thieme Verlag KG'
